$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.701.91'
$ws.Range('E2').Value = '  +2.13%  '
$ws.Range('D3').Value = '2.582.78'
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.52'
$ws.Range('E5').Value = '  -0.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.19'
$ws.Range('E6').Value = '  -2.51%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  +0.55%  '
$ws.Range('D9').Value = '2.592.38'
$ws.Range('E9').Value = '  +0.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.55'
$ws.Range('E10').Value = '  -1.06%  '
$ws.Range('E11').Value = '  +0.07%  '
$ws.Range('E12').Value = '  +1.34%  '
$ws.Range('E13').Value = '  +3.04%  '
$ws.Range('D14').Value = '3.036.32'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('D15').Value = '58.695.87'
$ws.Range('E15').Value = '  +2.15%  '
$ws.Range('E16').Value = '  +1.25%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.590.89'
$ws.Range('E17').Value = '  +1.83%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000132'
$ws.Range('E18').Value = '  -0.46%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '338.26'
$ws.Range('E19').Value = '  +1.07%  '
$ws.Range('E20').Value = '  +0.22%  '
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.50'
$ws.Range('E22').Value = '  +4.55%  '
$ws.Range('E23').Value = '  +0.08%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.31'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.404'
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('E27').Value = '  +0.12%  '
$ws.Range('E28').Value = '  +1.17%  '
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('E30').Value = '  -3.87%  '
$ws.Range('E31').Value = '  -5.05%  '
$ws.Range('E32').Value = '  -0.43%  '
$ws.Range('E33').Value = '  +0.98%  '
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.97'
$ws.Range('E35').Value = '  -1.38%  '
$ws.Range('E36').Value = '  -1.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '36.48'
$ws.Range('E37').Value = '  +1.11%  '
$ws.Range('E38').Value = '  +1.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.827'
$ws.Range('E39').Value = '  -0.16%  '
$ws.Range('E40').Value = '  -3.46%  '
$ws.Range('E41').Value = '  -0.24%  '
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '274.48'
$ws.Range('E43').Value = '  +2.13%  '
$ws.Range('E44').Value = '  +1.04%  '
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0948'
$ws.Range('E46').Value = '  -0.89%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0519'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.47'
$ws.Range('E48').Value = '  -1.91%  '
$ws.Range('D49').Value = '1.975.52'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('E50').Value = '  +0.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.50'
$ws.Range('E51').Value = '  -2.16%  '
